$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2099644128113879
$ws.Range("C2").Value = 0.5409252669039146
$ws.Range("J2").Value = 0.01423487544483986
$ws.Range("P2").Value = 0.1494661921708185
$ws.Range("S2").Value = 0.08540925266903915
$ws.Range("B3").Value = 0.01923076923076923
$ws.Range("C3").Value = 0.01282051282051282
$ws.Range("J3").Value = 0.03846153846153846
$ws.Range("P3").Value = 0.7307692307692307
$ws.Range("S3").Value = 0.1987179487179487
$ws.Range("O4").Value = 0.02439024390243903
$ws.Range("P4").Value = 0.7317073170731707
$ws.Range("S4").Value = 0.2439024390243902
$ws.Range("B6").Value = 0.07614213197969544
$ws.Range("F6").Value = 0.09644670050761421
$ws.Range("J6").Value = 0.1878172588832487
$ws.Range("O6").Value = 0.02030456852791878
$ws.Range("Q6").Value = 0.2487309644670051
$ws.Range("R6").Value = 0.05076142131979695
$ws.Range("S6").Value = 0.3197969543147208
$ws.Range("B7").Value = 0.08749999999999999
$ws.Range("D7").Value = 0.01875
$ws.Range("F7").Value = 0.075
$ws.Range("J7").Value = 0.125
$ws.Range("O7").Value = 0.01875
$ws.Range("Q7").Value = 0.15625
$ws.Range("R7").Value = 0.075
$ws.Range("S7").Value = 0.44375
$ws.Range("B8").Value = 0.09927360774818401
$ws.Range("D8").Value = 0.01694915254237288
$ws.Range("E8").Value = 0.002421307506053269
$ws.Range("F8").Value = 0.05811138014527845
$ws.Range("J8").Value = 0.09200968523002422
$ws.Range("O8").Value = 0.02421307506053269
$ws.Range("Q8").Value = 0.1961259079903148
$ws.Range("R8").Value = 0.1041162227602906
$ws.Range("S8").Value = 0.4067796610169492
$ws.Range("B9").Value = 0.1073170731707317
$ws.Range("D9").Value = 0.01951219512195122
$ws.Range("F9").Value = 0.06829268292682927
$ws.Range("J9").Value = 0.1170731707317073
$ws.Range("O9").Value = 0.01951219512195122
$ws.Range("Q9").Value = 0.1707317073170732
$ws.Range("R9").Value = 0.1121951219512195
$ws.Range("S9").Value = 0.3853658536585366
$ws.Range("B10").Value = 0.1047381546134663
$ws.Range("D10").Value = 0.02244389027431421
$ws.Range("F10").Value = 0.06400665004156277
$ws.Range("J10").Value = 0.1230257689110557
$ws.Range("O10").Value = 0.02410640066500416
$ws.Range("Q10").Value = 0.2036575228595179
$ws.Range("R10").Value = 0.07564422277639235
$ws.Range("S10").Value = 0.3823773898586866
$ws.Range("G11").Value = 0.1666666666666667
$ws.Range("J11").Value = 0.0945945945945946
$ws.Range("K11").Value = 0.2072072072072072
$ws.Range("L11").Value = 0.5315315315315315
$ws.Range("G12").Value = 0.7480314960629921
$ws.Range("J12").Value = 0.1889763779527559
$ws.Range("L12").Value = 0.03149606299212598
$ws.Range("S12").Value = 0.03149606299212598
$ws.Range("G13").Value = 0.6470588235294118
$ws.Range("J13").Value = 0.3137254901960784
$ws.Range("S13").Value = 0.0392156862745098
$ws.Range("F15").Value = 0.01941747572815534
$ws.Range("H15").Value = 0.1601941747572816
$ws.Range("I15").Value = 0.06310679611650485
$ws.Range("J15").Value = 0.354368932038835
$ws.Range("K15").Value = 0.04368932038834952
$ws.Range("M15").Value = 0.02427184466019417
$ws.Range("O15").Value = 0.07766990291262135
$ws.Range("S15").Value = 0.2572815533980582
$ws.Range("F16").Value = 0.02272727272727273
$ws.Range("H16").Value = 0.1590909090909091
$ws.Range("I16").Value = 0.1079545454545455
$ws.Range("J16").Value = 0.4204545454545455
$ws.Range("K16").Value = 0.07954545454545454
$ws.Range("M16").Value = 0.02840909090909091
$ws.Range("O16").Value = 0.07386363636363637
$ws.Range("S16").Value = 0.1079545454545455
$ws.Range("F17").Value = 0.0115473441108545
$ws.Range("H17").Value = 0.187066974595843
$ws.Range("I17").Value = 0.09930715935334873
$ws.Range("J17").Value = 0.4387990762124711
$ws.Range("K17").Value = 0.08314087759815242
$ws.Range("M17").Value = 0.02540415704387991
$ws.Range("N17").Value = 0.002309468822170901
$ws.Range("O17").Value = 0.04387990762124711
$ws.Range("S17").Value = 0.1085450346420323
$ws.Range("F18").Value = 0.005524861878453038
$ws.Range("H18").Value = 0.1712707182320442
$ws.Range("I18").Value = 0.09944751381215469
$ws.Range("J18").Value = 0.4972375690607735
$ws.Range("K18").Value = 0.0718232044198895
$ws.Range("M18").Value = 0.01657458563535912
$ws.Range("O18").Value = 0.03867403314917127
$ws.Range("S18").Value = 0.09944751381215469
$ws.Range("F19").Value = 0.01353637901861252
$ws.Range("H19").Value = 0.2072758037225042
$ws.Range("I19").Value = 0.09475465313028765
$ws.Range("J19").Value = 0.3866328257191201
$ws.Range("K19").Value = 0.08629441624365482
$ws.Range("M19").Value = 0.02622673434856176
$ws.Range("N19").Value = 0.0008460236886632825
$ws.Range("O19").Value = 0.06260575296108291
$ws.Range("S19").Value = 0.1218274111675127
